$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G (old G "d=7" shifts to H, old H "d=10" shifts to I)
$ws.Range("G1").EntireColumn.Insert()

# New header for the inserted column (match the header style used by the rest of row 1)
$ws.Range("G1").Value = "d=6"
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").HorizontalAlignment = -4108
$ws.Range("G1").VerticalAlignment = -4160
$ws.Range("G1").Borders.LineStyle = 1

# New values for the inserted column
$ws.Range("G2").Value = 97.89462345967142
$ws.Range("G3").Value = 98.04023545005992
$ws.Range("G4").Value = 97.96153372330222
$ws.Range("G5").Value = 97.90876057513898
$ws.Range("G6").Value = 97.95005638297364
